$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Details" text for the Test 2b / 14 GB GeoTiff rows: the
# benchmark was re-run on a 240 GB machine (was previously 64 GB).
$ws.Range("F142:F161").Value = "Read and decompress 14 GB GeoTiff from shared disk (norway_kartverket_10m_dtm_utm_z33.tif). 240 GB machine"

# Fill in the newly-collected "New" hub results for Test 2b (rows 152-161).
$ws.Range("E152").Value = 191.6
$ws.Range("E153").Value = 183.4
$ws.Range("E154").Value = 184.1
$ws.Range("E155").Value = 183.3
$ws.Range("E156").Value = 183.5
$ws.Range("E157").Value = 184
$ws.Range("E158").Value = 183.9
$ws.Range("E159").Value = 184.2
$ws.Range("E160").Value = 184.8
$ws.Range("E161").Value = 184.8

# Move the frozen-pane view / selection down to the newly added rows.
$ws.Application.ActiveWindow.ScrollRow = 134
$ws.Range("E162").Select()
